$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the text content of C2 and E2 (remove the "la identificacion..." and
# "El estado de cuenta..." shared-string values), while leaving I2/J2's
# styling intact but removing their string content too.
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Update the selection shown when the workbook is opened.
$ws.Range("B2:U4").Select()
